# #5: insurance, claim, debt, investment done
# Normalizes the "保險" (insurance) and "債務" (debt) sheets onto the
# common schema used across the other property sheets: replace the
# free-text "insurance period" / mortgage-detail column with a fixed
# "species" literal and append the shared trailer columns
# (property_category/category/date/legislator_name/legislator_id/
#  source_file/index).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "保險" (insurance) -> sheet9
# ---------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Header row: columns C/D/E get renamed to the normalized field names,
# and F:K are brand-new trailer headers.
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# Give the new header cells the same look (bold/border/center) as the
# rest of row 1.
$wsIns.Range("B1").Copy()
$wsIns.Range("F1:K1").PasteSpecial(-4122)

# Data rows 2-19: the old free-text "insurance period" column (E) is
# replaced by the fixed category literal "insurance"; F:J are constant
# across every row; K mirrors the existing index already in column A.
$lastRow = 19
for ($r = 2; $r -le $lastRow; $r++) {
    $wsIns.Range("E$r").Value = "insurance"
    $wsIns.Range("F$r").Value = "normal"
    $wsIns.Range("G$r").Value = "2012-04-27"
    $wsIns.Range("H$r").Value = "林滄敏"
    $wsIns.Range("I$r").Value = 1338
    $wsIns.Range("J$r").Value = "tmp9bfb1"
    $wsIns.Range("K$r").Value = $wsIns.Range("A$r").Value
}

$wsIns.Range("B2").Copy()
$wsIns.Range("F2:K$lastRow").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "債務" (debt) -> sheet10
# ---------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("債務")

# Header row: B/C/D/E/F/G get renamed to the normalized field names,
# and H:N are brand-new trailer headers.
$wsDebt.Range("B1").Value = "species"
$wsDebt.Range("C1").Value = "debtor"
$wsDebt.Range("D1").Value = "owner"
$wsDebt.Range("E1").Value = "total"
$wsDebt.Range("F1").Value = "register_date"
$wsDebt.Range("G1").Value = "register_reason"
$wsDebt.Range("H1").Value = "property_category"
$wsDebt.Range("I1").Value = "category"
$wsDebt.Range("J1").Value = "date"
$wsDebt.Range("K1").Value = "legislator_name"
$wsDebt.Range("L1").Value = "legislator_id"
$wsDebt.Range("M1").Value = "source_file"
$wsDebt.Range("N1").Value = "index"

$wsDebt.Range("B1").Copy()
$wsDebt.Range("H1:N1").PasteSpecial(-4122)

# Data row 2: old mortgage-detail column B keeps its literal value, but
# a new fixed "debt" category literal is introduced at H, and I:M are
# the same constant trailer values used on every sheet; N mirrors the
# existing index already in column A.
$wsDebt.Range("H2").Value = "debt"
$wsDebt.Range("I2").Value = "normal"
$wsDebt.Range("J2").Value = "2012-04-27"
$wsDebt.Range("K2").Value = "林滄敏"
$wsDebt.Range("L2").Value = 1338
$wsDebt.Range("M2").Value = "tmp9bfb1"
$wsDebt.Range("N2").Value = $wsDebt.Range("A2").Value

$wsDebt.Range("B2").Copy()
$wsDebt.Range("H2:N2").PasteSpecial(-4122)
